$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 38.5
$ws.Range("I9").Value = 30
$ws.Range("J9").Value = 55.5
$ws.Range("K9").Value = 30
$ws.Range("L9").Value = 55.5
$ws.Range("M9").Value = 139
$ws.Range("N9").Value = -393.5
$ws.Range("H19").Value = 446.22223
$ws.Range("I19").Value = 221
$ws.Range("K19").Value = 221
$ws.Range("M19").Value = -46
$ws.Range("H112").Value = 6037634.5
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 7094116
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 21282348
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -21284564
$ws.Range("H129").Value = 1250.1538
$ws.Range("I129").Value = 597
$ws.Range("J129").Value = 1304.5834
$ws.Range("K129").Value = 1791
$ws.Range("L129").Value = 3913.7502
$ws.Range("M129").Value = 3209
$ws.Range("N129").Value = -13913.7502
$ws.Range("H132").Value = 229257.03
$ws.Range("I132").Value = 253682.25
$ws.Range("J132").Value = 33855.332
$ws.Range("K132").Value = 761046.75
$ws.Range("L132").Value = 101565.996
$ws.Range("M132").Value = -758516.75
$ws.Range("N132").Value = -106625.996
$ws.Range("H135").Value = 1158.8375
$ws.Range("I135").Value = 1012.60657
$ws.Range("K135").Value = 9113.459130000001
$ws.Range("M135").Value = -6578.459130000001
$ws.Range("H137").Value = 20409216
$ws.Range("I137").Value = 29412532
$ws.Range("J137").Value = 1697.4
$ws.Range("K137").Value = 88237596
$ws.Range("L137").Value = 5092.200000000001
$ws.Range("M137").Value = -88235046
$ws.Range("N137").Value = -10192.2
$ws.Range("H138").Value = 8973283
$ws.Range("I138").Value = 1062062
$ws.Range("K138").Value = 3186186
$ws.Range("M138").Value = -3181046
$ws.Range("H141").Value = 1757.4667
$ws.Range("I141").Value = 1069.0308
$ws.Range("K141").Value = 3207.0924
$ws.Range("M141").Value = 1972.9076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 45609.305
$ws.Range("I2").Value = 69260.734
$ws.Range("K2").Value = 69260.734
$ws.Range("M2").Value = -69147.734
$ws.Range("H31").Value = 4433.3335
$ws.Range("I31").Value = 4433.3335
$ws.Range("K31").Value = 4433.3335
$ws.Range("M31").Value = -4139.3335
$ws.Range("H32").Value = 14272.149
$ws.Range("I32").Value = 2539.9578
$ws.Range("J32").Value = 66333.75
$ws.Range("K32").Value = 2539.9578
$ws.Range("L32").Value = 66333.75
$ws.Range("M32").Value = -2252.9578
$ws.Range("N32").Value = -66907.75
$ws.Range("H61").Value = 1312.58
$ws.Range("I61").Value = 1230.4468
$ws.Range("K61").Value = 1230.4468
$ws.Range("M61").Value = -1018.4468
$ws.Range("H74").Value = 3151.3594
$ws.Range("I74").Value = 1058.1063
$ws.Range("K74").Value = 1058.1063
$ws.Range("M74").Value = -184.1062999999999
$ws.Range("H77").Value = 3151.3594
$ws.Range("I77").Value = 1058.1063
$ws.Range("K77").Value = 5290.531499999999
$ws.Range("M77").Value = -922.5314999999991
$ws.Range("H116").Value = 45609.305
$ws.Range("I116").Value = 69260.734
$ws.Range("K116").Value = 69260.734
$ws.Range("M116").Value = -66966.734
$ws.Range("H132").Value = 1933.9783
$ws.Range("I132").Value = 1848.1282
$ws.Range("J132").Value = 2412.2856
$ws.Range("K132").Value = 5544.3846
$ws.Range("L132").Value = 7236.8568
$ws.Range("M132").Value = -3014.3846
$ws.Range("N132").Value = -12296.8568
$ws.Range("H136").Value = 1312.58
$ws.Range("I136").Value = 1230.4468
$ws.Range("K136").Value = 3691.3404
$ws.Range("M136").Value = -1141.3404

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 45609.305
$ws.Range("I3").Value = 69260.734
$ws.Range("K3").Value = 69260.734
$ws.Range("M3").Value = -69146.734
$ws.Range("H81").Value = 29336
$ws.Range("J81").Value = 29336
$ws.Range("L81").Value = 29336
$ws.Range("N81").Value = -31458
$ws.Range("H84").Value = 29336
$ws.Range("J84").Value = 29336
$ws.Range("L84").Value = 88008
$ws.Range("N84").Value = -98616
$ws.Range("H134").Value = 11906012
$ws.Range("I134").Value = 12821663
$ws.Range("J134").Value = 2551
$ws.Range("K134").Value = 38464989
$ws.Range("L134").Value = 7653
$ws.Range("M134").Value = -38462454
$ws.Range("N134").Value = -12723

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 333336670
$ws.Range("J4").Value = 333336670
$ws.Range("L4").Value = 333336670
$ws.Range("N4").Value = -333336894
$ws.Range("H58").Value = 1688.0952
$ws.Range("I58").Value = 928.7742
$ws.Range("J58").Value = 3828
$ws.Range("K58").Value = 928.7742
$ws.Range("L58").Value = 3828
$ws.Range("M58").Value = -725.7742
$ws.Range("N58").Value = -4234
$ws.Range("H99").Value = 15625725
$ws.Range("I99").Value = 20833800
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 20833800
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -20832302
$ws.Range("N99").Value = -4496
$ws.Range("H126").Value = 15625725
$ws.Range("I126").Value = 20833800
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 62501400
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -62498930
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 1612.7119
$ws.Range("I132").Value = 1155.1372
$ws.Range("J132").Value = 4529.75
$ws.Range("K132").Value = 3465.411599999999
$ws.Range("L132").Value = 13589.25
$ws.Range("M132").Value = -935.4115999999995
$ws.Range("N132").Value = -18649.25
$ws.Range("H134").Value = 1854.5968
$ws.Range("I134").Value = 1225.8704
$ws.Range("K134").Value = 3677.6112
$ws.Range("M134").Value = -1142.6112
$ws.Range("H136").Value = 1688.0952
$ws.Range("I136").Value = 928.7742
$ws.Range("J136").Value = 3828
$ws.Range("K136").Value = 2786.3226
$ws.Range("L136").Value = 11484
$ws.Range("M136").Value = -236.3226
$ws.Range("N136").Value = -16584

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3005.3635
$ws.Range("I139").Value = 2823.5264
$ws.Range("K139").Value = 8470.5792
$ws.Range("M139").Value = -3330.5792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2618.4
$ws.Range("I102").Value = 2382.1538
$ws.Range("K102").Value = 2382.1538
$ws.Range("M102").Value = -760.1538
$ws.Range("H132").Value = 2336.0476
$ws.Range("I132").Value = 2061.5
$ws.Range("J132").Value = 3983.3333
$ws.Range("K132").Value = 6184.5
$ws.Range("L132").Value = 11949.9999
$ws.Range("M132").Value = -3654.5
$ws.Range("N132").Value = -17009.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2487.2727
$ws.Range("I46").Value = 2100
$ws.Range("K46").Value = 2100
$ws.Range("M46").Value = -1912
$ws.Range("H122").Value = 3382.4443
$ws.Range("I122").Value = 1996
$ws.Range("J122").Value = 3778.5715
$ws.Range("K122").Value = 5988
$ws.Range("L122").Value = 11335.7145
$ws.Range("M122").Value = -3538
$ws.Range("N122").Value = -16235.7145
$ws.Range("H132").Value = 4439.222
$ws.Range("I132").Value = 4767.8047
$ws.Range("J132").Value = 3402.923
$ws.Range("K132").Value = 14303.4141
$ws.Range("L132").Value = 10208.769
$ws.Range("M132").Value = -11773.4141
$ws.Range("N132").Value = -15268.769
$ws.Range("H136").Value = 2550.6528
$ws.Range("I136").Value = 1556.2239
$ws.Range("J136").Value = 15876
$ws.Range("K136").Value = 4668.6717
$ws.Range("L136").Value = 47628
$ws.Range("M136").Value = -2118.6717
$ws.Range("N136").Value = -52728

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 50027500
$ws.Range("J2").Value = 60001.5
$ws.Range("L2").Value = 60001.5
$ws.Range("N2").Value = -60225.5
$ws.Range("H58").Value = 12547
$ws.Range("J58").Value = 12547
$ws.Range("L58").Value = 12547
$ws.Range("N58").Value = -13163
$ws.Range("H62").Value = 7711725
$ws.Range("I62").Value = 16692083
$ws.Range("K62").Value = 16692083
$ws.Range("M62").Value = -16691459
$ws.Range("H65").Value = 7711725
$ws.Range("I65").Value = 16692083
$ws.Range("K65").Value = 83460415
$ws.Range("M65").Value = -83457295
$ws.Range("H122").Value = 54376.95
$ws.Range("I122").Value = 92463.73
$ws.Range("K122").Value = 277391.19
$ws.Range("M122").Value = -274941.19
$ws.Range("H132").Value = 11630268
$ws.Range("I132").Value = 16131401
$ws.Range("J132").Value = 2341.6667
$ws.Range("K132").Value = 48394203
$ws.Range("L132").Value = 7025.000100000001
$ws.Range("M132").Value = -48391673
$ws.Range("N132").Value = -12085.0001
$ws.Range("H136").Value = 18542.543
$ws.Range("I136").Value = 19445.74
$ws.Range("K136").Value = 58337.22
$ws.Range("M136").Value = -55787.22
